$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'27.254.06"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Formula = "'1.564.49"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Formula = "'210.66"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").Formula = "'1.787.32"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Formula = "'1.576.38"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Formula = "'27.213.34"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Formula = "'217.93"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Formula = "'0.0₃0703"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").Formula = "'1.95"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Formula = "'151.61"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").Formula = "'6.62"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").Formula = "'15.02"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Formula = "'3.18"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Formula = "'1.455.79"
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("D35").Formula = "'1.12"
$ws.Range("E35").Value = "  +5.20%  "
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D41").Formula = "'0.814"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").Formula = "'64.49"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Formula = "'1.699.10"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Formula = "'85.93"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Formula = "'0.0523"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  -1.66%  "
